# Custom Error Code Details - cleanup spell/grammar proofing marks and
# populate the trailing blank row with the new 4041 "User Data  not found"
# entry.
#
# NOTE: simply assigning Range.Text (or Find/Replace) across a run that is
# wrapped by a leading <w:proofErr/> pair can leave an orphaned proofErr
# element behind in the saved XML, because the proofing marks are siblings
# of the run rather than part of it. Deleting the *whole* cell range first
# (which clears the paragraph completely, proofErr marks included) and then
# inserting the final text back in gives a clean single <w:r> with no
# leftover <w:proofErr/> markers - matching how Word normally cleans these
# up once the flagged text is edited.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellPlainText($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Delete()
    $cell.Range.InsertBefore($text)
}

function Find-RowByCellText($table, $needle) {
    for ($i = 1; $i -le $table.Rows.Count; $i++) {
        $cellText = $table.Cell($i, 2).Range.Text
        if ($cellText.Contains($needle)) {
            return $i
        }
    }
    return -1
}

function Find-EmptyDataRow($table) {
    for ($i = 2; $i -le $table.Rows.Count; $i++) {
        $c1 = $table.Cell($i, 1).Range.Text
        $c2 = $table.Cell($i, 2).Range.Text
        # Cell.Range.Text always ends with the paragraph mark + cell mark
        # (two characters), so strip those before checking for emptiness.
        $c1Content = $c1.Substring(0, $c1.Length - 2)
        $c2Content = $c2.Substring(0, $c2.Length - 2)
        if ($c1Content -eq "" -and $c2Content -eq "") {
            return $i
        }
    }
    return -1
}

# 1) "Sql Grammar Exception" - drop the spellStart/spellEnd around "Sql"
$row = Find-RowByCellText $t "Sql"
Set-CellPlainText $t $row 2 "Sql Grammar Exception"

# 2) "Network or driver issue or db is temporarily unavailable" - drop the
#    spellStart/spellEnd around "db"
$row = Find-RowByCellText $t "Network or driver issue or"
Set-CellPlainText $t $row 2 "Network or driver issue or db is temporarily unavailable"

# 3) "Sorry we could not found your profile information" - drop the
#    gramStart/gramEnd around "found"
$row = Find-RowByCellText $t "Sorry we could not"
Set-CellPlainText $t $row 2 "Sorry we could not found your profile information"

# 4) "UserType is Wrong" - drop the spellStart/spellEnd around "UserType"
$row = Find-RowByCellText $t "is Wrong"
Set-CellPlainText $t $row 2 "UserType is Wrong"

# 5) Populate the trailing blank row with the new error code entry.
$row = Find-EmptyDataRow $t
Set-CellPlainText $t $row 1 "4041"
Set-CellPlainText $t $row 2 "User Data  not found"
